# EIA Table 2.4.C monthly update: 2006-October 2016 -> 2006-November 2016
# Adds the "November" monthly row into the "Year to Date" block, shifts the
# trailing summary rows down by one, refreshes their totals, and updates the
# "Rolling 12 Months Ending in ..." label + totals, plus the report title.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row 53 (between "October" row 52 and the "Year to Date"
#    section header, currently row 53) to hold the new "November" data
#    point, pushing every row from 53 on down by one.
# ---------------------------------------------------------------------
$ws.Rows(53).Insert()

# Copy the number formatting (styles) from an existing month-data row (row 39,
# "November" in the "Year 2016" block) onto the newly inserted row so the new
# row matches the surrounding month rows exactly.
$fmtSrc = $ws.Range("A39:F39")
$newRow = $ws.Range("A53:F53")
$fmtSrc.Copy()
$newRow.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new "November" row in the "Year to Date" section.
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 777672
$ws.Range("C53").Value = 333264
$ws.Range("D53").Value = 329725
$ws.Range("E53").Value = 9130
$ws.Range("F53").Value = 105553

# ---------------------------------------------------------------------
# 2) Refresh the "Year to Date" annual totals (now rows 55-57, were 54-56)
#    to include the new November figures.
# ---------------------------------------------------------------------
$ws.Range("B55").Value = 8661019
$ws.Range("C55").Value = 3595746
$ws.Range("D55").Value = 3911467
$ws.Range("E55").Value = 108709
$ws.Range("F55").Value = 1045098

$ws.Range("B56").Value = 10063086
$ws.Range("C56").Value = 4359214
$ws.Range("D56").Value = 4485550
$ws.Range("E56").Value = 106611
$ws.Range("F56").Value = 1111712

$ws.Range("B57").Value = 10575283
$ws.Range("C57").Value = 4712735
$ws.Range("D57").Value = 4606771
$ws.Range("E57").Value = 107008
$ws.Range("F57").Value = 1148769

# ---------------------------------------------------------------------
# 3) Update the "Rolling 12 Months Ending in ..." section (now rows 58-60,
#    were 57-59): relabel the header and refresh its totals.
# ---------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

$ws.Range("B59").Value = 10811599
$ws.Range("C59").Value = 4663402
$ws.Range("D59").Value = 4820131
$ws.Range("E59").Value = 116493
$ws.Range("F59").Value = 1211573

$ws.Range("B60").Value = 11463871
$ws.Range("C60").Value = 5106836
$ws.Range("D60").Value = 4981277
$ws.Range("E60").Value = 116777
$ws.Range("F60").Value = 1258981

# ---------------------------------------------------------------------
# 4) Update the report title/subtitle to reflect the new "through November"
#    coverage period.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Million Cubic Feet)"

Write-Output "done"
